# Apply "Updated with protocol links" changes to Sheet1.
# Schema (row 1 headers): A=id, B=parent_id, C=position, D=step_name,
# E=step_desc, F=organism, G=sop_url, H=duration_in_days,
# I=part_of_service, J=comments, K=success_rate, L=created, M=last_modified

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 98: SEVA/METABRICK, Moclo, Gibson pathway engineering @Wageningen ---
$ws.Range("E98").ClearContents()
$ws.Range("F98").ClearContents()
$ws.Range("G98").ClearContents()
$ws.Range("I98").ClearContents()
$ws.Range("M98").Value = "2019-10-11 11:23:41"

# --- Row 99: becomes the "Vector, promoter selection, ..." task ---
$ws.Range("A99").Value = 67
$ws.Range("B99").Value = 64
$ws.Range("C99").Value = 2
$ws.Range("D99").Value = "***** Vector, promoter selection, plasmids, locus of integration, restriction sites, etc (cloning strategy) #SopNeeded"
$ws.Range("E99").ClearContents()
$ws.Range("F99").ClearContents()
$ws.Range("G99").ClearContents()
$ws.Range("I99").ClearContents()
$ws.Range("M99").Value = "2019-10-11 11:23:41"

# --- Row 100: new task "Select primers for vector assembly" ---
$ws.Range("A100").Value = 400
$ws.Range("B100").Value = 67
$ws.Range("C100").Value = 1
$ws.Range("D100").Value = "****** Select primers for vector assembly"
$ws.Range("E100").Value = "Primers required for construction of the vector"
$ws.Range("F100").ClearContents()
$ws.Range("G100").Value = "https://hub.ibisba.eu/sops/58"
$ws.Range("I100").ClearContents()
$ws.Range("L100").Value = 43749.47356481482
$ws.Range("M100").Value = "2019-10-11 11:26:03"

# --- Row 101: Combinatorial assembly strategy #SopNeeded ---
$ws.Range("E101").ClearContents()
$ws.Range("F101").ClearContents()
$ws.Range("G101").ClearContents()
$ws.Range("I101").ClearContents()
$ws.Range("M101").Value = "2019-10-11 11:23:41"

# --- Row 136: Prepare liquid cultivation medium for transformation ---
$ws.Range("G136").Value = "https://hub.ibisba.eu/sops/56"
$ws.Range("M136").Value = "2019-10-11 11:28:03"

# --- Row 137: Prepare solid cultivation medium for transformation ---
$ws.Range("G137").Value = "https://hub.ibisba.eu/sops/55"
$ws.Range("M137").Value = "2019-10-11 11:32:47"

# --- Row 140: Make glycerol stock ---
$ws.Range("G140").Value = "https://hub.ibisba.eu/sops/54"
$ws.Range("M140").Value = "2019-10-11 11:35:20"

# --- Row 144: In vitro assembly of the DNA construct ---
$ws.Range("F144").ClearContents()
$ws.Range("G144").Value = "https://hub.ibisba.eu/sops/52"
$ws.Range("I144").ClearContents()
$ws.Range("M144").Value = "2019-10-11 11:41:13"

# --- Row 146: Prepare cultivation media for transformation ---
$ws.Range("F146").ClearContents()
$ws.Range("G146").Value = "https://hub.ibisba.eu/sops/56"
$ws.Range("I146").ClearContents()
$ws.Range("M146").Value = "2019-10-11 11:30:51"

# --- Row 157: Make glycerol stock ---
$ws.Range("D157").Value = "****** Make glycerol stock"
$ws.Range("F157").ClearContents()
$ws.Range("G157").Value = "https://hub.ibisba.eu/sops/54, https://hub.ibisba.eu/sops/53"
$ws.Range("I157").ClearContents()
$ws.Range("M157").Value = "2019-10-11 13:05:31"

# --- Row 158: Add label/barcode ---
$ws.Range("D158").Value = "****** Add label/barcode"
$ws.Range("F158").ClearContents()
$ws.Range("I158").ClearContents()
$ws.Range("M158").Value = "2019-10-11 13:05:45"

# --- Row 160: In vitro assembly of the integration/expression vector ---
$ws.Range("F160").ClearContents()
$ws.Range("G160").Value = "https://hub.ibisba.eu/sops/57"
$ws.Range("I160").ClearContents()
$ws.Range("M160").Value = "2019-10-11 11:41:31"

# --- Row 173: Make glycerol stock ---
$ws.Range("F173").ClearContents()
$ws.Range("G173").Value = "https://hub.ibisba.eu/sops/54, https://hub.ibisba.eu/sops/53"
$ws.Range("I173").ClearContents()
$ws.Range("M173").Value = "2019-10-11 11:37:35"

# --- Row 179: Prepare genetic material for transformation ---
$ws.Range("F179").ClearContents()
$ws.Range("G179").Value = "https://hub.ibisba.eu/sops/51"
$ws.Range("I179").ClearContents()
$ws.Range("M179").Value = "2019-10-11 11:38:22"

# --- Row 196: Recombinant clone screening and selection ---
$ws.Range("D196").Value = "**** Recombinant clone screening and selection (depends on transformation implementation)"
$ws.Range("E196").ClearContents()
$ws.Range("F196").ClearContents()
$ws.Range("G196").ClearContents()
$ws.Range("I196").ClearContents()
$ws.Range("M196").Value = "2019-10-11 13:07:27"

# --- Row 199: Make glycerol stock #SopNeeded ---
$ws.Range("F199").ClearContents()
$ws.Range("G199").Value = "https://hub.ibisba.eu/sops/54, https://hub.ibisba.eu/sops/53"
$ws.Range("I199").ClearContents()
$ws.Range("M199").Value = "2019-10-11 11:39:11"
